# Fruta / hortaliza, semanal
# Insert a new weekly record as row 80, shifting the existing rows 80-137 down
# to 81-138 (dimension grows from A1:R137 to A1:R138).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 80; this shifts rows 80..137 down
# to 81..138, preserving their values/formatting (including the date style
# on column D).
$ws.Rows("80").Insert()

# Populate the newly inserted row 80 with the new record's data.
$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "Macroferia Regional de Talca"
$ws.Range("C80").Value = "Maule"
$ws.Range("D80").Value = 45161
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = 100112013
$ws.Range("G80").Value = "Alcachofa"
$ws.Range("H80").Value = "Madrigal"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 200
$ws.Range("K80").Value = 13000
$ws.Range("L80").Value = 13000
$ws.Range("M80").Value = 13000
$ws.Range("N80").Value = "$/caja 40 unidades"
$ws.Range("O80").Value = "Provincia del Elquí"
$ws.Range("P80").Value = 325
$ws.Range("Q80").Value = 40
$ws.Range("R80").Value = "Hortaliza"
